$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6811249999999999
$ws.Range("M2").Value = 0.9979466666666666
$ws.Range("N2").Value = 2.99384
$ws.Range("O2").Value = 0.4487505885107415
$ws.Range("P2").Value = 0.4487505885107414
$ws.Range("Q2").Value = 0.6797264233333332
$ws.Range("R2").Value = 6.117537809999999
$ws.Range("S2").Value = 0.4487505885107415
$ws.Range("T2").Value = 0.4487505885107414
$ws.Range("G3").Value = 0.6811249999999999
$ws.Range("O3").Value = 0.2623432118199488
$ws.Range("P3").Value = 0.2623432118199487
$ws.Range("Q3").Value = 0.3973735469583333
$ws.Range("S3").Value = 0.2623432118199488
$ws.Range("T3").Value = 0.2623432118199487
$ws.Range("G4").Value = 0.6811249999999999
$ws.Range("M4").Value = 0.4584083333333333
$ws.Range("N4").Value = 1.375225
$ws.Range("O4").Value = 0.206134271732853
$ws.Range("P4").Value = 0.2061342717328529
$ws.Range("Q4").Value = 0.3122333760416666
$ws.Range("R4").Value = 2.810100384375
$ws.Range("S4").Value = 0.206134271732853
$ws.Range("T4").Value = 0.2061342717328529
$ws.Range("G5").Value = 0.6811249999999999
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.184071
$ws.Range("N5").Value = 0.5522130000000001
$ws.Range("O5").Value = 0.08277192793645689
$ws.Range("P5").Value = 0.08277192793645688
$ws.Range("Q5").Value = 0.125375359875
$ws.Range("R5").Value = 1.128378238875
$ws.Range("S5").Value = 0.08277192793645689
$ws.Range("T5").Value = 0.08277192793645688
